# Apply updates to the "strategy" worksheet:
#  - Update set_voltage (column G) values for rows 5 and 8-18
#  - Update the active cell selection to G11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (set_voltage) values
$ws.Range("G5").Value = 54.6
$ws.Range("G8").Value = 51
$ws.Range("G9").Value = 51
$ws.Range("G10").Value = 51
$ws.Range("G11").Value = 51
$ws.Range("G12").Value = 51
$ws.Range("G13").Value = 51
$ws.Range("G14").Value = 51
$ws.Range("G15").Value = 51
$ws.Range("G16").Value = 51
$ws.Range("G17").Value = 51
$ws.Range("G18").Value = 51

# Move/set the active selection to G11 (matches the updated sheetView selection)
$ws.Range("G11").Select()
